$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-18 05:42:11"
$wsZh.Range("G2").Value = "2016-02-18 05:42:59"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-18 05:42:27"
$wsDe.Range("G2").Value = "2016-02-18 05:43:29"
